$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New metric values (same across all data rows 2-26), per the updated training run.
$values = @{
    "B" = 0.9999824846980448
    "C" = 0.9990072066285755
    "D" = 0.9999703103901542
    "E" = 0.9999998042879644
    "F" = 0.9999924385405996
    "G" = 0.00001634976554226276
    "H" = 0.0009267290336323482
    "I" = 0.00001563167196290704
    "J" = 0.00000007975901297912329
    "K" = 0.00000785571548794308
    "L" = 0.0002547310818618769
    "M" = 0.004043484331892824
    "N" = 1.000016814689877
    "O" = 0.004215623831067973
    "P" = 120.0425940012932
    "Q" = 179.767509419835
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
